# 2nd version, added "add window"
#
# - "Novembro" sheet: the trailing formatted-but-empty row (row 26, which only
#   carried a styled blank cell in column B) is deleted, shifting it up to
#   become the new row 25 and shrinking the used range by one row.
# - "Dezembro" sheet: the placeholder/junk entry rows (A4:D5 - "fghdf",
#   "asdfsdf", "saf", "asd", the stray "30/11/2023" text dates, etc.) are
#   removed, while the adjoining summary cells in columns G/H stay untouched.

$wb = $excel.ActiveWorkbook

$wsNovembro = $wb.Worksheets("Novembro")
$wsNovembro.Rows(25).Delete()

$wsDezembro = $wb.Worksheets("Dezembro")
$wsDezembro.Range("A4:D5").Clear()
